# Update cryptos list data (prices / volume %) and reorder three rows
# (Mantle / BabyDogeCoin / EnergySwap) per the Aug 1 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.917.95"
$ws.Range("E2").Value = "  -1.72%  "
$ws.Range("D3").Value = "'1.833.36"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "'244.98"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'0.6908"
$ws.Range("E6").Value = "  -2.12%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "'0.07682"
$ws.Range("E8").Value = "  -2.65%  "
$ws.Range("D9").Value = "'0.3052"
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("D10").Value = "'23.52"
$ws.Range("E10").Value = "  -4.15%  "
$ws.Range("D11").Value = "'0.07814"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").Value = "'1.829.26"
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").Value = "'5.072"
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("D14").Value = "'90.53"
$ws.Range("E14").Value = "  -3.23%  "
$ws.Range("D15").Value = "'0.6802"
$ws.Range("E15").Value = "  -2.96%  "
$ws.Range("D16").Value = "'6.435"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").Value = "'0.000008346"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "'28.923.03"
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("D19").Value = "'243.20"
$ws.Range("E19").Value = "  -4.13%  "
$ws.Range("D20").Value = "'2.082.72"
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("D21").Value = "'12.68"
$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "'7.479"
$ws.Range("E23").Value = "  -1.94%  "
$ws.Range("D24").Value = "'0.9997"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").Value = "'0.1467"
$ws.Range("E25").Value = "  -5.90%  "
$ws.Range("D26").Value = "'161.68"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'8.798"
$ws.Range("E27").Value = "  -2.20%  "
$ws.Range("E28").Value = "  -3.41%  "
$ws.Range("D29").Value = "'1.555"
$ws.Range("E29").Value = "  +3.00%  "
$ws.Range("D30").Value = "'4.214"
$ws.Range("E30").Value = "  -2.69%  "
$ws.Range("D31").Value = "'4.154"
$ws.Range("E31").Value = "  -2.45%  "
$ws.Range("D32").Value = "'1.177"
$ws.Range("E32").Value = "  -2.56%  "
$ws.Range("D33").Value = "'0.05123"
$ws.Range("E33").Value = "  -3.20%  "
$ws.Range("D34").Value = "'0.7652"
$ws.Range("E34").Value = "  +2.17%  "
$ws.Range("D35").Value = "'1.847"
$ws.Range("E35").Value = "  -2.53%  "
$ws.Range("D36").Value = "'1.146"
$ws.Range("E36").Value = "  -2.21%  "
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").Value = "'0.01843"
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("D39").Value = "'1.231.57"
$ws.Range("E39").Value = "  -3.70%  "
$ws.Range("D40").Value = "'2.692"
$ws.Range("E40").Value = "  -2.79%  "
$ws.Range("D41").Value = "'0.9206"
$ws.Range("E41").Value = "  +2.84%  "
$ws.Range("D42").Value = "'108.42"
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("D43").Value = "'5.843"
$ws.Range("E43").Value = "  -3.21%  "
$ws.Range("D44").Value = "'0.9997"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'9.559"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "'1.983.61"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.5169"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000122"
$ws.Range("E48").Value = "  -4.78%  "
$ws.Range("D49").Value = "'63.88"
$ws.Range("E49").Value = "  -10.23%  "
$ws.Range("D50").Value = "'1.745"
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("D51").Value = "'6.927"
$ws.Range("E51").Value = "  -1.74%  "

# The apostrophe prefix above also stamps each D-column cell with a
# "quote prefix" style; reset back to the Normal style so the cell keeps
# its original (unstyled) text value with no leftover formatting flag.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"

